$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, 29380, "Diogo Viana", "TI", "Outros", 6, 45098, 10404.22),
  @(3, 75412, "Pedro Henrique Mendes", "Operações", "Doença", 6, 45093, 9850.82),
  @(4, 60658, "Maria Fernanda Moura", "Atendimento ao Cliente", "Outros", 4, 45083, 9435.58),
  @(5, 19416, "Maria Julia Ramos", "TI", "Problemas pessoais", 5, 45081, 3130.98),
  @(6, 37749, "Natália Gomes", "Engenharia", "Doença", 5, 45094, 9674.62),
  @(7, 41751, "Cecília Duarte", "Marketing", "Problemas pessoais", 5, 45085, 9582.58),
  @(8, 43462, "Olivia Silva", "Vendas", "Doença", 1, 45090, 3736.63),
  @(9, 15243, "Bruna Rezende", "Financeiro", "Outros", 1, 45085, 3263.23),
  @(10, 79587, "Fernando da Mota", "TI", "Outros", 2, 45089, 8839.09),
  @(11, 91696, "Juliana Azevedo", "Operações", "Viagem de negócios", 8, 45093, 8928.5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
